$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Cells.Item(2, 4)
$c.NumberFormat = "@"
$c.Value2 = '26.703.34'
$c.Style = "Normal"
$c = $ws.Cells.Item(3, 4)
$c.NumberFormat = "@"
$c.Value2 = '1.599.70'
$c.Style = "Normal"
$ws.Cells.Item(3, 5).Value = '  +0.35%  '
$ws.Cells.Item(4, 5).Value = '  +0.15%  '
$c = $ws.Cells.Item(5, 4)
$c.NumberFormat = "@"
$c.Value2 = '211.38'
$c.Style = "Normal"
$ws.Cells.Item(5, 5).Value = '  +0.25%  '
$ws.Cells.Item(6, 5).Value = '  -0.35%  '
$ws.Cells.Item(7, 5).Value = '  +0.12%  '
$ws.Cells.Item(8, 5).Value = '  +0.50%  '
$ws.Cells.Item(9, 5).Value = '  +1.18%  '
$ws.Cells.Item(10, 5).Value = '  +0.67%  '
$ws.Cells.Item(11, 5).Value = '  +0.87%  '
$c = $ws.Cells.Item(12, 4)
$c.NumberFormat = "@"
$c.Value2 = '1.824.63'
$c.Style = "Normal"
$ws.Cells.Item(12, 5).Value = '  +0.35%  '
$c = $ws.Cells.Item(13, 4)
$c.NumberFormat = "@"
$c.Value2 = '1.601.39'
$c.Style = "Normal"
$ws.Cells.Item(13, 5).Value = '  -0.49%  '
$c = $ws.Cells.Item(14, 4)
$c.NumberFormat = "@"
$c.Value2 = '4.04'
$c.Style = "Normal"
$ws.Cells.Item(14, 5).Value = '  +0.71%  '
$ws.Cells.Item(15, 5).Value = '  +0.79%  '
$c = $ws.Cells.Item(16, 4)
$c.NumberFormat = "@"
$c.Value2 = '65.29'
$c.Style = "Normal"
$ws.Cells.Item(16, 5).Value = '  +1.38%  '
$c = $ws.Cells.Item(17, 4)
$c.NumberFormat = "@"
$c.Value2 = '26.684.08'
$c.Style = "Normal"
$ws.Cells.Item(17, 5).Value = '  +0.38%  '
$c = $ws.Cells.Item(18, 4)
$c.NumberFormat = "@"
$c.Value2 = '0.0₃0756'
$c.Style = "Normal"
$ws.Cells.Item(18, 5).Value = '  +3.63%  '
$c = $ws.Cells.Item(19, 4)
$c.NumberFormat = "@"
$c.Value2 = '209.83'
$c.Style = "Normal"
$ws.Cells.Item(19, 5).Value = '  +1.02%  '
$ws.Cells.Item(20, 2).Value = 'Dai'
$ws.Cells.Item(20, 3).Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$c = $ws.Cells.Item(20, 4)
$c.NumberFormat = "@"
$c.Value2 = '1.00'
$c.Style = "Normal"
$ws.Cells.Item(20, 5).Value = '  +0.16%  '
$ws.Cells.Item(21, 2).Value = 'Chainlink'
$ws.Cells.Item(21, 3).Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$c = $ws.Cells.Item(21, 4)
$c.NumberFormat = "@"
$c.Value2 = '7.19'
$c.Style = "Normal"
$ws.Cells.Item(21, 5).Value = '  +4.46%  '
$ws.Cells.Item(22, 5).Value = '  +0.84%  '
$ws.Cells.Item(23, 5).Value = '  +0.15%  '
$c = $ws.Cells.Item(25, 4)
$c.NumberFormat = "@"
$c.Value2 = '143.10'
$c.Style = "Normal"
$ws.Cells.Item(25, 5).Value = '  -1.45%  '
$ws.Cells.Item(26, 5).Value = '  +0.13%  '
$ws.Cells.Item(27, 5).Value = '  -0.02%  '
$c = $ws.Cells.Item(28, 4)
$c.NumberFormat = "@"
$c.Value2 = '0.115'
$c.Style = "Normal"
$ws.Cells.Item(28, 5).Value = '  +0.43%  '
$c = $ws.Cells.Item(29, 4)
$c.NumberFormat = "@"
$c.Value2 = '15.34'
$c.Style = "Normal"
$ws.Cells.Item(29, 5).Value = '  +0.81%  '
$c = $ws.Cells.Item(30, 4)
$c.NumberFormat = "@"
$c.Value2 = '0.0517'
$c.Style = "Normal"
$ws.Cells.Item(30, 5).Value = '  +2.74%  '
$ws.Cells.Item(31, 5).Value = '  +0.09%  '
$ws.Cells.Item(32, 5).Value = '  +1.12%  '
$ws.Cells.Item(33, 5).Value = '  +1.83%  '
$c = $ws.Cells.Item(34, 4)
$c.NumberFormat = "@"
$c.Value2 = '1.290.86'
$c.Style = "Normal"
$c = $ws.Cells.Item(35, 4)
$c.NumberFormat = "@"
$c.Value2 = '0.619'
$c.Style = "Normal"
$ws.Cells.Item(35, 5).Value = '  -4.82%  '
$ws.Cells.Item(36, 5).Value = '  +0.97%  '
$ws.Cells.Item(37, 5).Value = '  +0.60%  '
$ws.Cells.Item(38, 5).Value = '  -0.03%  '
$ws.Cells.Item(39, 5).Value = '  +15.71%  '
$ws.Cells.Item(40, 5).Value = '  -1.81%  '
$ws.Cells.Item(41, 5).Value = '  -0.66%  '
$c = $ws.Cells.Item(42, 4)
$c.NumberFormat = "@"
$c.Value2 = '0.786'
$c.Style = "Normal"
$ws.Cells.Item(42, 5).Value = '  +0.17%  '
$ws.Cells.Item(43, 5).Value = '  -0.52%  '
$ws.Cells.Item(44, 5).Value = '  -0.80%  '
$c = $ws.Cells.Item(45, 4)
$c.NumberFormat = "@"
$c.Value2 = '1.736.63'
$c.Style = "Normal"
$c = $ws.Cells.Item(46, 4)
$c.NumberFormat = "@"
$c.Value2 = '91.02'
$c.Style = "Normal"
$ws.Cells.Item(46, 5).Value = '  +1.73%  '
$c = $ws.Cells.Item(47, 4)
$c.NumberFormat = "@"
$c.Value2 = '1.57'
$c.Style = "Normal"
$ws.Cells.Item(47, 5).Value = '  -0.71%  '
$c = $ws.Cells.Item(48, 4)
$c.NumberFormat = "@"
$c.Value2 = '0.101'
$c.Style = "Normal"
$ws.Cells.Item(48, 5).Value = '  -1.36%  '
$ws.Cells.Item(49, 5).Value = '  +0.77%  '
$ws.Cells.Item(50, 5).Value = '  +0.14%  '
$c = $ws.Cells.Item(51, 4)
$c.NumberFormat = "@"
$c.Value2 = '7.35'
$c.Style = "Normal"
$ws.Cells.Item(51, 5).Value = '  -1.09%  '
